# Adds the "http://purl.org/dc/terms/isVersionOf" column to every sheet of
# the workbook, populated with the AWV asset URI that corresponds to each
# row's own asset id. On sheet 1 ("onderdeel#HeeftAanzicht") the new column
# is simply appended after the last used column; on the other four sheets
# it is inserted right before the existing "isActief" / RelatieObject.*
# column(s), shifting those columns one position to the right.

$wb = $excel.ActiveWorkbook

$isVersionOf = "http://purl.org/dc/terms/isVersionOf"

# --- Sheet 1: onderdeel#HeeftAanzicht --------------------------------------
# Dimension A1:J2 -> A1:K2 (plain append, no shift needed).
$ws1 = $wb.Worksheets.Item("onderdeel#HeeftAanzicht")
$ws1.Range("K1").Value = $isVersionOf
$ws1.Range("K2").Value = "https://data.awvvlaanderen.be/id/asset/HeeftAanzicht_-_opstelling_01_-_aanzicht_01"

# --- Sheet 2: onderdeel#HoortBij -------------------------------------------
# Dimension A1:H2 -> A1:I2. New column inserted at G, pushing the
# RelatieObject.bron/doel columns from G,H to H,I.
$ws2 = $wb.Worksheets.Item("onderdeel#HoortBij")
$ws2.Columns.Item(7).Insert()
$ws2.Range("G1").Value = $isVersionOf
$ws2.Range("G2").Value = "https://data.awvvlaanderen.be/id/asset/HoortBij_-_bord_01_-_aanzicht_01"

# --- Sheet 3: installatie#AanzichtVerkeersbordopstelling -------------------
# Dimension A1:I2 -> A1:J2. New column inserted at H, pushing
# isActief/toestand from H,I to I,J.
$ws3 = $wb.Worksheets.Item("installatie#AanzichtVerkeersbordopstelling")
$ws3.Columns.Item(8).Insert()
$ws3.Range("H1").Value = $isVersionOf
$ws3.Range("H2").Value = "https://data.awvvlaanderen.be/id/asset/aanzicht_01"

# --- Sheet 4: onderdeel#RetroreflecterendVerkeersbord ----------------------
# Dimension A1:J3 -> A1:K3. New column inserted at H, pushing
# isActief/opstelhoogte/toestand from H,I,J to I,J,K. Both data rows
# (bord_01) get the same isVersionOf value.
$ws4 = $wb.Worksheets.Item("onderdeel#RetroreflecterendVerkeersbord")
$ws4.Columns.Item(8).Insert()
$ws4.Range("H1").Value = $isVersionOf
$ws4.Range("H2").Value = "https://data.awvvlaanderen.be/id/asset/bord_01"
$ws4.Range("H3").Value = "https://data.awvvlaanderen.be/id/asset/bord_01"

# --- Sheet 5: installatie#Verkeersbordopstelling ----------------------------
# Dimension A1:I2 -> A1:J2. New column inserted at H, pushing
# isActief/toestand from H,I to I,J.
$ws5 = $wb.Worksheets.Item("installatie#Verkeersbordopstelling")
$ws5.Columns.Item(8).Insert()
$ws5.Range("H1").Value = $isVersionOf
$ws5.Range("H2").Value = "https://data.awvvlaanderen.be/id/asset/opstelling_01"
